$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update filenames and Param1 values for all data rows, and reorder the
# 8bit / Fusion blocks (8bit now precedes Fusion).
$ws.Cells.Item(2, 1).Value = "Zelda--param1-00.68.dac"
$ws.Cells.Item(2, 2).Value = "Zelda"
$ws.Cells.Item(2, 3).Value = 0.68
$ws.Cells.Item(3, 1).Value = "Zelda--param1-00.53.dac"
$ws.Cells.Item(3, 2).Value = "Zelda"
$ws.Cells.Item(3, 3).Value = 0.53
$ws.Cells.Item(4, 1).Value = "Zelda--param1-01.00.dac"
$ws.Cells.Item(4, 2).Value = "Zelda"
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(5, 1).Value = "Zelda--param1-00.60.dac"
$ws.Cells.Item(5, 2).Value = "Zelda"
$ws.Cells.Item(5, 3).Value = 0.6
$ws.Cells.Item(6, 1).Value = "Zelda--param1-00.32.dac"
$ws.Cells.Item(6, 2).Value = "Zelda"
$ws.Cells.Item(6, 3).Value = 0.32
$ws.Cells.Item(7, 1).Value = "lofi--param1-00.49.dac"
$ws.Cells.Item(7, 2).Value = "lofi"
$ws.Cells.Item(7, 3).Value = 0.49
$ws.Cells.Item(8, 1).Value = "lofi--param1-00.60.dac"
$ws.Cells.Item(8, 2).Value = "lofi"
$ws.Cells.Item(8, 3).Value = 0.6
$ws.Cells.Item(9, 1).Value = "lofi--param1-00.88.dac"
$ws.Cells.Item(9, 2).Value = "lofi"
$ws.Cells.Item(9, 3).Value = 0.88
$ws.Cells.Item(10, 1).Value = "lofi--param1-00.44.dac"
$ws.Cells.Item(10, 2).Value = "lofi"
$ws.Cells.Item(10, 3).Value = 0.44
$ws.Cells.Item(11, 1).Value = "lofi--param1-00.25.dac"
$ws.Cells.Item(11, 2).Value = "lofi"
$ws.Cells.Item(11, 3).Value = 0.25
$ws.Cells.Item(12, 1).Value = "8bit--param1-00.90.dac"
$ws.Cells.Item(12, 2).Value = "8bit"
$ws.Cells.Item(12, 3).Value = 0.9
$ws.Cells.Item(13, 1).Value = "8bit--param1-00.02.dac"
$ws.Cells.Item(13, 2).Value = "8bit"
$ws.Cells.Item(13, 3).Value = 0.02
$ws.Cells.Item(14, 1).Value = "8bit--param1-00.35.dac"
$ws.Cells.Item(14, 2).Value = "8bit"
$ws.Cells.Item(14, 3).Value = 0.35
$ws.Cells.Item(15, 1).Value = "8bit--param1-00.13.dac"
$ws.Cells.Item(15, 2).Value = "8bit"
$ws.Cells.Item(15, 3).Value = 0.13
$ws.Cells.Item(16, 1).Value = "8bit--param1-00.85.dac"
$ws.Cells.Item(16, 2).Value = "8bit"
$ws.Cells.Item(16, 3).Value = 0.85
$ws.Cells.Item(17, 1).Value = "Fusion--param1-00.82.dac"
$ws.Cells.Item(17, 2).Value = "Fusion"
$ws.Cells.Item(17, 3).Value = 0.82
$ws.Cells.Item(18, 1).Value = "Fusion--param1-00.42.dac"
$ws.Cells.Item(18, 2).Value = "Fusion"
$ws.Cells.Item(18, 3).Value = 0.42
$ws.Cells.Item(19, 1).Value = "Fusion--param1-00.39.dac"
$ws.Cells.Item(19, 2).Value = "Fusion"
$ws.Cells.Item(19, 3).Value = 0.39
$ws.Cells.Item(20, 1).Value = "Fusion--param1-00.34.dac"
$ws.Cells.Item(20, 2).Value = "Fusion"
$ws.Cells.Item(20, 3).Value = 0.34
$ws.Cells.Item(21, 1).Value = "Fusion--param1-00.11.dac"
$ws.Cells.Item(21, 2).Value = "Fusion"
$ws.Cells.Item(21, 3).Value = 0.11
$ws.Cells.Item(22, 1).Value = "duduk--param1-00.66.dac"
$ws.Cells.Item(22, 2).Value = "duduk"
$ws.Cells.Item(22, 3).Value = 0.66
$ws.Cells.Item(23, 1).Value = "duduk--param1-00.03.dac"
$ws.Cells.Item(23, 2).Value = "duduk"
$ws.Cells.Item(23, 3).Value = 0.03
$ws.Cells.Item(24, 1).Value = "duduk--param1-00.38.dac"
$ws.Cells.Item(24, 2).Value = "duduk"
$ws.Cells.Item(24, 3).Value = 0.38
$ws.Cells.Item(25, 1).Value = "duduk--param1-00.11.dac"
$ws.Cells.Item(25, 2).Value = "duduk"
$ws.Cells.Item(25, 3).Value = 0.11
$ws.Cells.Item(26, 1).Value = "duduk--param1-00.30.dac"
$ws.Cells.Item(26, 2).Value = "duduk"
$ws.Cells.Item(26, 3).Value = 0.3

# Remove the obsolete Param2 column (D) entirely.
$ws.Range("D1:D26").EntireColumn.Delete()
